$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2235
$ws1.Range("F3").Value = 102
$ws1.Range("F4").Value = 13518
$ws1.Range("F7").Value = 532
$ws1.Range("F11").Value = 13851
$ws1.Range("F12").Value = 14610
$ws1.Range("F26").Value = 5611
$ws1.Range("F29").Value = 5370
$ws1.Range("F30").Value = 39
$ws1.Range("F32").Value = 186

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2235
$ws4.Range("F3").Value = 102
$ws4.Range("F4").Value = 13518
$ws4.Range("F8").Value = 532
$ws4.Range("F12").Value = 13851
$ws4.Range("F13").Value = 14610
$ws4.Range("F27").Value = 5611
$ws4.Range("F30").Value = 5370
$ws4.Range("F31").Value = 39
$ws4.Range("F33").Value = 186
